$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Variables": insert a new row for "urb_area_id_preg" right after
# the existing "urb_area_id_0" row, and reword the description of
# "urb_area_id_0" to clarify it refers to "at birth".
# ---------------------------------------------------------------------
$wsVar = $wb.Worksheets.Item("Variables")

$wsVar.Rows.Item(105).Insert()
$wsVar.Range("A105:D105").Style = "Normal"
$wsVar.Range("A105").Value = "urb_area_id_preg"
$wsVar.Range("B105").Value = "integer"
$wsVar.Range("D105").Value = "unique identifier for the urban area at pregnancy (for the cohorts with children moving from one area to another at different time periods) "
$wsVar.Range("D104").Value = "unique identifier for the urban area at birth (for the cohorts with children moving from one area to another at different time periods) "

# Re-apply the autofilter so its range grows to match the new last row.
$wsVar.AutoFilterMode = $false
$wsVar.Range("A1:D219").AutoFilter()

# Keep the workbook-level _FilterDatabase defined name in sync.
$fdVar = $wb.Names.Item("Variables!_FilterDatabase")
$fdVar.RefersTo = "=Variables!`$A`$1:`$D`$219"

# ---------------------------------------------------------------------
# Sheet "Categories": insert 4 new rows right after the existing
# "urb_area_id_0" category rows, duplicating the same city list but tied
# to the new "urb_area_id_preg" variable.
# ---------------------------------------------------------------------
$wsCat = $wb.Worksheets.Item("Categories")

$wsCat.Rows("79:82").Insert()
$wsCat.Range("A79:D82").Style = "Normal"

$wsCat.Range("A79").Value = "urb_area_id_preg"
$wsCat.Range("B79").Value = 1401
$wsCat.Range("C79").Value = $false
$wsCat.Range("D79").Value = "Copenhagen"

$wsCat.Range("A80").Value = "urb_area_id_preg"
$wsCat.Range("B80").Value = 1402
$wsCat.Range("C80").Value = $false
$wsCat.Range("D80").Value = "Odense"

$wsCat.Range("A81").Value = "urb_area_id_preg"
$wsCat.Range("B81").Value = 1403
$wsCat.Range("C81").Value = $false
$wsCat.Range("D81").Value = "Aarhus"

$wsCat.Range("A82").Value = "urb_area_id_preg"
$wsCat.Range("B82").Value = 1404
$wsCat.Range("C82").Value = $false
$wsCat.Range("D82").Value = "Aalborg"

$wsCat.Range("A3").Select()

# Leave the original sheet ("Variables") focused/active with A3 selected,
# matching the saved view state.
$wsVar.Activate()
$wsVar.Range("A3").Select()
